$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("R40" label in the rule table) is changed from the text "R40"
# to the text "1". Force a Text number format first so Excel stores the
# new value as a shared string (text) rather than silently coercing the
# numeric-looking "1" into a Number cell, matching the original edit
# (the cell keeps its string/text type in the OOXML, just a new shared
# string value).
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
